$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet1): bump "想去人数" (F column) counts for several rows.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(5,6).Value  = 15486
$wsExpo.Cells.Item(8,6).Value  = 698
$wsExpo.Cells.Item(9,6).Value  = 15375
$wsExpo.Cells.Item(11,6).Value = 8976
$wsExpo.Cells.Item(18,6).Value = 195
$wsExpo.Cells.Item(20,6).Value = 44
$wsExpo.Cells.Item(25,6).Value = 1105
$wsExpo.Cells.Item(28,6).Value = 78
$wsExpo.Cells.Item(34,6).Value = 246
$wsExpo.Cells.Item(37,6).Value = 114
$wsExpo.Cells.Item(38,6).Value = 5499

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet2): append a new event row (row 4) for the Vienna Royal
# Philharmonic New Year concert.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)

# Copy row 3's formatting (bold/bordered index cell) down into row 4 first so
# the new index cell (column A) keeps the same style as the existing ones.
$wsShow.Cells.Item(3,1).Copy($wsShow.Cells.Item(4,1))

$wsShow.Cells.Item(4,1).Value = 3
$wsShow.Cells.Item(4,2).NumberFormat = "@"
$wsShow.Cells.Item(4,2).Value = "2024-12-22"
$wsShow.Cells.Item(4,3).Value = "苏州·维也纳皇家交响乐团2025新年音乐会"
$wsShow.Cells.Item(4,4).Value = "东苑路1号公共文化中心内 苏州保利大剧院"
$wsShow.Cells.Item(4,5).Value = "2024.12.22 19:30-12.22 21:30"
$wsShow.Cells.Item(4,6).Value = 0
$wsShow.Cells.Item(4,7).Value = 280
$wsShow.Cells.Item(4,8).Value = "https://show.bilibili.com/platform/detail.html?id=92817"
$wsShow.Cells.Item(4,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/cCTiFEpg1727155421223.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet4): same F-column bumps as "展览" (rows shifted by
# the extra rows this combined sheet already carries), plus insertion of the
# same new concert event between the existing row 40 and row 41 (so the old
# "星部落动漫嘉年华" row slides from row 41 down to row 42).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(5,6).Value  = 15486
$wsAll.Cells.Item(8,6).Value  = 698
$wsAll.Cells.Item(9,6).Value  = 15375
$wsAll.Cells.Item(11,6).Value = 8976
$wsAll.Cells.Item(18,6).Value = 195
$wsAll.Cells.Item(20,6).Value = 44
$wsAll.Cells.Item(25,6).Value = 1105
$wsAll.Cells.Item(28,6).Value = 78
$wsAll.Cells.Item(36,6).Value = 246
$wsAll.Cells.Item(39,6).Value = 114
$wsAll.Cells.Item(40,6).Value = 5499

# Shift the existing row 41 ("星部落动漫嘉年华") down to row 42, carrying its
# formatting with it, then renumber its index cell.
$wsAll.Range("A41:I41").Copy($wsAll.Range("A42:I42"))
$wsAll.Cells.Item(42,1).Value = 41

# Overwrite row 41 in place with the new concert event (same data as the new
# row added to "演出"), keeping the existing bold/bordered index-cell style.
$wsAll.Cells.Item(41,1).Value = 40
$wsAll.Cells.Item(41,2).NumberFormat = "@"
$wsAll.Cells.Item(41,2).Value = "2024-12-22"
$wsAll.Cells.Item(41,3).Value = "苏州·维也纳皇家交响乐团2025新年音乐会"
$wsAll.Cells.Item(41,4).Value = "东苑路1号公共文化中心内 苏州保利大剧院"
$wsAll.Cells.Item(41,5).Value = "2024.12.22 19:30-12.22 21:30"
$wsAll.Cells.Item(41,6).Value = 0
$wsAll.Cells.Item(41,7).Value = 280
$wsAll.Cells.Item(41,8).Value = "https://show.bilibili.com/platform/detail.html?id=92817"
$wsAll.Cells.Item(41,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/cCTiFEpg1727155421223.jpeg"
